# Rename the three "rework" sheets to the final "Test Case N" names and
# refresh the view state (selection / zoom) that Excel recorded for them.
#
# TestReport -> Test Case 1
# Rework1    -> Test Case 2
# Rework2    -> Test Case 3

$wb = $excel.ActiveWorkbook

$wsTestCase1 = $wb.Worksheets.Item("TestReport")
$wsTestCase2 = $wb.Worksheets.Item("Rework1")
$wsTestCase3 = $wb.Worksheets.Item("Rework2")

$wsTestCase1.Name = "Test Case 1"
$wsTestCase2.Name = "Test Case 2"
$wsTestCase3.Name = "Test Case 3"

# The sheet-scoped "Print_Titles" defined name keeps the OLD sheet name
# baked into its formula until the page-setup print titles are re-applied,
# so nudge each sheet's PageSetup to force it to pick up the new name.
$wsTestCase1.PageSetup.PrintTitleRows = "`$8:`$8"
$wsTestCase2.PageSetup.PrintTitleRows = "`$8:`$8"
$wsTestCase3.PageSetup.PrintTitleRows = "`$8:`$8"

# "Test Case 1" (formerly TestReport) had the user's selection sitting on
# C9; move it to B10.
[void]$wsTestCase1.Select()
[void]$wsTestCase1.Range("B10").Select()

# "Test Case 2" (formerly Rework1) was zoomed out from 85% to 70%.
[void]$wsTestCase2.Select()
$excel.ActiveWindow.Zoom = 70

# Leave "Test Case 1" as the active/selected tab, matching the workbook's
# recorded activeTab.
[void]$wsTestCase1.Select()
